$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 175, shifting existing rows 175-289 down to 176-290.
$ws.Rows.Item(175).Insert()

# Fill the newly inserted row 175 with the new weekly record.
$ws.Cells.Item(175, 1).Value = 5
$ws.Cells.Item(175, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(175, 3).Value = "Maule"
$ws.Cells.Item(175, 4).Value = 44719
$ws.Cells.Item(175, 5).Value = 7
$ws.Cells.Item(175, 6).Value = 100112006
$ws.Cells.Item(175, 7).Value = "Repollo"
$ws.Cells.Item(175, 8).Value = "Crespo record"
$ws.Cells.Item(175, 9).Value = "Primera"
$ws.Cells.Item(175, 10).Value = 3000
$ws.Cells.Item(175, 11).Value = 1200
$ws.Cells.Item(175, 12).Value = 1200
$ws.Cells.Item(175, 13).Value = 1200
$ws.Cells.Item(175, 14).Value = "`$/unidad"
$ws.Cells.Item(175, 15).Value = "Región del Maule"
$ws.Cells.Item(175, 16).Value = 1200
$ws.Cells.Item(175, 17).Value = 1
$ws.Cells.Item(175, 18).Value = "Hortaliza"
